$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.760400000000001
$ws.Range("C2").Value = 0.7485000000000001

$ws.Range("B3").Value = 3.0618
$ws.Range("C3").Value = 0.6493

$ws.Range("B4").Value = 10.969
$ws.Range("C4").Value = 0.6844

$ws.Range("B5").Value = 3.1462
$ws.Range("C5").Value = 0.5704

$ws.Range("B6").Value = 5.8306
$ws.Range("C6").Value = 0.8656

$ws.Range("B7").Value = 1.114
$ws.Range("C7").Value = 0.6481

$ws.Range("B8").Value = 3.5094
$ws.Range("C8").Value = 0.8032

$ws.Range("B9").Value = 2.8994
$ws.Range("C9").Value = 0.6581

$ws.Range("B10").Value = 0.5479000000000001
$ws.Range("C10").Value = 0.5454
